$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, forcing text storage so
# numeric-looking strings (trailing zeros, exact decimals, leading zeros)
# are preserved exactly like the workbook's existing inline-string cells,
# then restore the Normal style so no extraneous formatting is introduced.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "62.734.66"
$ws.Range("E2").Value = "  -2.55%  "

# Row 3
$ws.Range("D3").Value = "3.377.42"
$ws.Range("E3").Value = "  -3.70%  "

# Row 4
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
Set-TextValue "D5" "571.95"
$ws.Range("E5").Value = "  -3.29%  "

# Row 6
Set-TextValue "D6" "125.58"
$ws.Range("E6").Value = "  -6.66%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "3.376.64"
$ws.Range("E8").Value = "  -3.70%  "

# Row 9
Set-TextValue "D9" "0.477"
$ws.Range("E9").Value = "  -2.25%  "

# Row 10
Set-TextValue "D10" "7.29"
$ws.Range("E10").Value = "  -4.29%  "

# Row 11
Set-TextValue "D11" "0.119"
$ws.Range("E11").Value = "  -4.83%  "

# Row 12
Set-TextValue "D12" "0.376"
$ws.Range("E12").Value = "  -3.36%  "

# Row 13
$ws.Range("D13").Value = "3.950.00"
$ws.Range("E13").Value = "  -3.78%  "

# Row 14
$ws.Range("E14").Value = "  -0.85%  "

# Row 15
$ws.Range("D15").Value = "3.371.20"
$ws.Range("E15").Value = "  -3.96%  "

# Row 16
Set-TextValue "D16" "0.0000170"
$ws.Range("E16").Value = "  -6.06%  "

# Row 17
$ws.Range("D17").Value = "62.751.22"
$ws.Range("E17").Value = "  -2.49%  "

# Row 18
Set-TextValue "D18" "24.38"
$ws.Range("E18").Value = "  -5.38%  "

# Row 19
Set-TextValue "D19" "9.23"
$ws.Range("E19").Value = "  -7.88%  "

# Row 20
Set-TextValue "D20" "5.62"
$ws.Range("E20").Value = "  -2.41%  "

# Row 21
Set-TextValue "D21" "13.07"
$ws.Range("E21").Value = "  -4.16%  "

# Row 22
Set-TextValue "D22" "371.61"
$ws.Range("E22").Value = "  -5.17%  "

# Row 23
Set-TextValue "D23" "0.552"
$ws.Range("E23").Value = "  -4.86%  "

# Row 24
$ws.Range("D24").Value = "3.509.00"
$ws.Range("E24").Value = "  -3.78%  "

# Row 25
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
Set-TextValue "D26" "71.44"
$ws.Range("E26").Value = "  -4.07%  "

# Row 27
Set-TextValue "D27" "0.0000105"
$ws.Range("E27").Value = "  -10.59%  "

# Row 28
Set-TextValue "D28" "0.995"
$ws.Range("E28").Value = "  -0.50%  "

# Row 29
Set-TextValue "D29" "6.95"
$ws.Range("E29").Value = "  -6.41%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "7.80"
$ws.Range("E30").Value = "  -5.14%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.10"
$ws.Range("E31").Value = "  -7.57%  "

# Row 32
$ws.Range("E32").Value = "  -0.02%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D33" "1.38"
$ws.Range("E33").Value = "  -6.31%  "

# Row 34
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.405.50"
$ws.Range("E34").Value = "  -3.63%  "

# Row 35
Set-TextValue "D35" "0.148"
$ws.Range("E35").Value = "  -5.91%  "

# Row 36
Set-TextValue "D36" "22.66"
$ws.Range("E36").Value = "  -3.05%  "

# Row 37
Set-TextValue "D37" "5.40"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
Set-TextValue "D38" "166.20"
$ws.Range("E38").Value = "  -0.35%  "

# Row 39
Set-TextValue "D39" "6.62"
$ws.Range("E39").Value = "  -4.94%  "

# Row 40
Set-TextValue "D40" "1.47"
$ws.Range("E40").Value = "  -5.81%  "

# Row 41
Set-TextValue "D41" "0.0753"
$ws.Range("E41").Value = "  -4.46%  "

# Row 42
$ws.Range("E42").Value = "  -0.08%  "

# Row 43
Set-TextValue "D43" "41.72"
$ws.Range("E43").Value = "  -0.84%  "

# Row 44
Set-TextValue "D44" "0.764"
$ws.Range("E44").Value = "  -5.80%  "

# Row 45
Set-TextValue "D45" "4.21"
$ws.Range("E45").Value = "  -5.31%  "

# Row 46
Set-TextValue "D46" "1.54"
$ws.Range("E46").Value = "  -7.36%  "

# Row 47
Set-TextValue "D47" "22.49"
$ws.Range("E47").Value = "  -10.05%  "

# Row 48
Set-TextValue "D48" "1.09"
$ws.Range("E48").Value = "  -7.75%  "

# Row 49
Set-TextValue "D49" "6.58"
$ws.Range("E49").Value = "  -3.24%  "

# Row 50
$ws.Range("D50").Value = "2.230.63"
$ws.Range("E50").Value = "  -6.07%  "

# Row 51
Set-TextValue "D51" "0.830"
$ws.Range("E51").Value = "  -8.89%  "
